# Fall 2021 schedule: push quiz availability from lecture day to the
# following class, and relabel the quizzes as "open" on that day.
# (commit: "update syllabus push lecture 7")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Quiz 3: move from row 11 (M 9/27) to row 12 (W 9/29), label "Quiz 3 open"
$ws.Range("D11").ClearContents()
$ws.Range("D12").Value = "Quiz 3 open"

# Quiz 4: move from row 13 (M 10/4) to row 14 (W 10/6), label "Quiz 4 open"
$ws.Range("D13").ClearContents()
$ws.Range("D14").Value = "Quiz 4 open"

# Quiz 5: move from row 17 (M 10/18) to row 18 (W 10/20), label "Quiz 5 open"
$ws.Range("D17").ClearContents()
$ws.Range("D18").Value = "Quiz 5 open"

# Quiz 6: move from row 19 (M 10/25) to row 20 (W 10/27), label "Quiz 6 open"
$ws.Range("D19").ClearContents()
$ws.Range("D20").Value = "Quiz 6 open"

# Quiz 7: move from row 21 (M 11/1) to row 22 (W 11/3), label "Quiz 7 open"
$ws.Range("D21").ClearContents()
$ws.Range("D22").Value = "Quiz 7 open"

# Quiz 8: move from row 23 (M 11/8) to row 24 (W 11/10), label "Quiz 8 open"
$ws.Range("D23").ClearContents()
$ws.Range("D24").Value = "Quiz 8 open"

# Update the saved view/selection to match where the edit left off.
$ws.Range("D25").Select() | Out-Null
